$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.728.81'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '3.504.14'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '180.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.11%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.505.79'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.594'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('E10').Value = '  +7.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.13'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.03%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').Value = '4.111.44'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.40'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +11.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.137'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.92%  '
$ws.Range('D16').Value = '67.722.85'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000178'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '3.514.01'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.31'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '396.58'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.539'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.71'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.178'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.35%  '
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('E32').Value = '  -0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.08'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.72'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '164.18'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.876'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('E41').Value = '  +6.96%  '
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.22%  '
$ws.Range('D44').Value = '2.861.77'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '26.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0730'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.56'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '42.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.16%  '
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '341.26'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.70%  '
